# Highlight quantitative/impact metrics (percentages, dollar amounts,
# large numbers) in bold + color (#2C3E50) across the resume body.
#
# Word's Font.Color COM property takes a BGR-packed integer (classic
# Windows COLORREF 0x00BBGGRR), while OOXML <w:color w:val="RRGGBB"/>
# stores RRGGBB -- so convert once and reuse.
$HighlightColor = 0x2C3E50
$BgrColor = (($HighlightColor -band 0xFF) * 0x10000) + ($HighlightColor -band 0xFF00) + (($HighlightColor -band 0xFF0000) / 0x10000)

$d = $word.ActiveDocument
$PM = [char]0xB1     # "±" plus/minus sign
$Bullet = [char]0x2022   # "•" bullet

# Bold + color the first occurrence of $needle found at/after $searchStart
# within the paragraph ending at $paraEnd. Returns the End position of the
# matched (and now highlighted) range so callers can chain further
# searches later in the same paragraph.
function HighlightNext($searchStart, $paraEnd, $needle) {
    $scope = $d.Range($searchStart, $paraEnd)
    $found = $scope.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find '$needle' in range $searchStart-$paraEnd"
    }
    $scope.Font.Bold = 1
    $scope.Font.Color = $BgrColor
    return $scope.End
}

# Locate the single paragraph whose text (sans trailing paragraph mark)
# equals $exactText exactly, and apply bold+color to each needle in
# $needles, in left-to-right order.
function HighlightMetricsInParagraph($exactText, $needles) {
    $count = $d.Paragraphs.Count
    $target = $null
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $range = $para.Range
        $t = $range.Text
        $t = $t.Substring(0, $t.Length - 1)
        if ($t -eq $exactText) {
            $target = $range
            break
        }
    }
    if ($target -eq $null) {
        throw "Could not find paragraph matching '$exactText'"
    }
    $cursor = $target.Start
    $paraEnd = $target.End
    foreach ($needle in $needles) {
        $cursor = HighlightNext $cursor $paraEnd $needle
    }
}

# 1) "...improving demographic classification accuracy from 23% to 64%"
HighlightMetricsInParagraph `
    ($Bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%") `
    @("23%", "64%")

# 2) "Achieved 87% ... standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
HighlightMetricsInParagraph `
    ($Bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $PM + "4.2% to " + $PM + "2.1%") `
    @("87%", "71%", ($PM + "4.2%"), ($PM + "2.1%"))

# 3) "Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
HighlightMetricsInParagraph `
    ($Bullet + " Wrote RFP and analyzed bids from 1,200 vendors for research platform development") `
    @("1,200")

# 4) "...became the $400M Polling Consortium Database ... now valued at $1B+"
HighlightMetricsInParagraph `
    ($Bullet + " Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+") `
    @("`$400M", "`$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
HighlightMetricsInParagraph `
    ($Bullet + " Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M") `
    @("73.5%", "`$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (Key Achievements)
HighlightMetricsInParagraph `
    ($Bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%") `
    @("87%", "71%")

Write-Output "Highlighted metrics in 6 paragraphs."
